$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Dice-sum table (columns I/J/K/M, rows 3-13):
#   I = sum of two dice (2..12)
#   J = probability of that sum
#   K = I * J         (contribution to the mean)
#   M = (I-mean)^2*J  (contribution to the variance)
# ---------------------------------------------------------------
for ($r = 3; $r -le 13; $r++) {
    $ws.Cells.Item($r, 9).Value = $r - 1   # column I: 2..12
}

$ws.Range("J3").Formula = "=1/36"
$ws.Range("J4").Formula = "=2/36"
$ws.Range("J5").Formula = "=3/36"
$ws.Range("J6").Formula = "=4/36"
$ws.Range("J7").Formula = "=5/36"
$ws.Range("J8").Formula = "=6/36"
$ws.Range("J9").Value = 0.1388888888888889
$ws.Range("J10").Formula = "=4/36"
$ws.Range("J11").Formula = "=3/36"
$ws.Range("J12").Formula = "=2/36"
$ws.Range("J13").Formula = "=1/36"
$ws.Range("J3:J13").NumberFormat = "0.00"

$ws.Range("K3").Formula = "=I3*J3"
$ws.Range("K4:K13").Formula = "=I4*J4"

$ws.Range("M3").Formula = "=((I3-`$K`$14)^2)*J3"
$ws.Range("M4:M13").Formula = "=((I4-`$K`$14)^2)*J4"

# ---------------------------------------------------------------
# Summary row 14/15: mean / var / std labels + totals.
# Written in this order so the new shared-string table gets
# "var", "std", "mean" (in that order).
# ---------------------------------------------------------------
$ws.Range("L14").Value = "var"
$ws.Range("L15").Value = "std"
$ws.Range("J14").Value = "mean"

$ws.Range("K14").Formula = "=SUM(K3:K13)"
$ws.Range("K14").NumberFormat = "0.00"

$ws.Range("M14").Formula = "=SUM(M3:M13)"
$ws.Range("M14").NumberFormat = "0.00"

$ws.Range("N14").Formula = "=M14*50"

$ws.Range("M15").Formula = "=M14^(1/2)"
$ws.Range("N15").Formula = "=N14^(1/2)"

# ---------------------------------------------------------------
# Roulette bet table (rows 20-23)
# ---------------------------------------------------------------
$ws.Range("E20").Value = 35
$ws.Range("F20").Formula = "=1/38"
$ws.Range("G20").Formula = "=E20*F20"
$ws.Range("I20").Formula = "=((E20-G22)^2)*F20"

$ws.Range("E21").Value = -1
$ws.Range("F21").Formula = "=37/38"
$ws.Range("G21").Formula = "=E21*F21"
$ws.Range("I21").Formula = "=((E21-G22)^2)*F21"

$ws.Range("G22").Formula = "=SUM(G20:G21)"
$ws.Range("I22").Formula = "=SUM(I20:I21)"
$ws.Range("J22").Formula = "=I22*100"

$ws.Range("G23").Formula = "=G22*100"
$ws.Range("I23").Formula = "=I22^(1/2)"
$ws.Range("J23").Formula = "=J22^(1/2)"

# ---------------------------------------------------------------
# Cosmetic view/column adjustments to mirror the authored workbook
# (widened to fit the new K/M numeric columns, ~10.57 chars wide)
# ---------------------------------------------------------------
$ws.Columns("K").ColumnWidth = 9.736979166666666
$ws.Columns("M").ColumnWidth = 9.736979166666666

$ws.Range("M24").Select()
